$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("1")

$data = @(
  @(2,  0.2607118305826748,  86.58224122722923,  -0.194251745467877),
  @(3,  0.1433229263445609,  86.30074735584871,   0.4760979816206624),
  @(4,  0.3465224987268718,  82.40498098273099,  -0.04617925570788317),
  @(5,  0.6606286740243559,  81.35155381000143,   1.086085773919071),
  @(6,  0.7072059114790494,  76.9053271022166,    1.275674131713703),
  @(7,  1.281452089395786,   67.92011851261924,   2.284045568691468),
  @(8,  1.757875618971923,   61.07412900997206,   3.924738434293051),
  @(9,  2.475190580137981,   49.72493395585144,   6.049749363973344),
  @(10, 3.735754587525885,   35.34485309594061,   8.530198545172324),
  @(11, 4.532281842628445,   22.68575716653009,  10.05161169375461),
  @(12, 5.148166713156415,   13.48379381638977,  12.63227643877123),
  @(13, 5.514710387128434,    6.446273358922222, 13.57815389603109),
  @(14, 5.890766284955653,    1.808955719152385, 13.81788816785999),
  @(15, 6.203789973888749,    2.069059665242438, 14.05552423271236),
  @(16, 6.151247439233878,    0.8700556024477172,14.00863237921378),
  @(17, 6.087556708946686,    0.56017547751703,  14.46836756712724),
  @(18, 6.064159626381498,    0.2280575701626303,14.35688141571183),
  @(19, 6.156098899253257,    0.5916137758930876,14.71925356196701),
  @(20, 5.972101842355398,   -0.3248492833437174,14.17440967383386),
  @(21, 6.298214609974625,   -1.034272526608184, 14.49957502304193)
)

foreach ($row in $data) {
  $r = $row[0]
  $ws.Cells.Item($r, 3).Value = $row[1]
  $ws.Cells.Item($r, 4).Value = $row[2]
  $ws.Cells.Item($r, 5).Value = $row[3]
}
